$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 10) for 2021, matching the structure of the
# preceding yearly rows (row 9 = 2020).
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 39.5081226894545
$ws.Range("C10").Value = 25.1231084154561
$ws.Range("D10").Value = 130.004872252279
$ws.Range("E10").Value = 31.3520249812064
$ws.Range("F10").Value = 60.7625801155657
$ws.Range("G10").Value = 13.9868880989956
$ws.Range("H10").Value = 13.2511200026315
$ws.Range("I10").Value = 39.8905678125791
$ws.Range("J10").Value = 57.7502398556484
$ws.Range("K10").Value = 106.11188888254
$ws.Range("L10").Value = 70.1736416228309
$ws.Range("M10").Value = 35.7761168680666
$ws.Range("N10").Value = 14.8602119724028
$ws.Range("O10").Value = 70.8253348807258
$ws.Range("P10").Value = 88.7216618728709
$ws.Range("Q10").Value = 25.6510579035764
$ws.Range("R10").Value = 26.9386217306364

# Match the style of column A's year-label cells (bold text with a thin
# border, as used by A2:A9) by copying the format from the cell directly
# above.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
